# Remove "rad" from a handful of radiometric-routine related strings in the
# document body (mirrors the blunt "rad" -> "" substring removal the author
# applied across the package).  Each Find/Replace below targets one exact,
# uniquely-occurring phrase so nothing else in the document is touched.

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Exact($findText, $replaceText) {
    $d.Content.Find.Execute(
        $findText,
        $true,
        $false,
        $false,
        $false,
        $false,
        $true,
        $wdFindContinue,
        $false,
        $replaceText,
        $wdReplaceAll
    ) | Out-Null
}

Replace-Exact "Converts to radiance 31.3 mWatts/[cm2-sr] (cf. 31.3)" "Converts to iance 31.3 mWatts/[cm2-sr] (cf. 31.3)"
Replace-Exact "Converts to corneal irradiance 29.9 uWatts/cm2 (cf. 29.9)" "Converts to corneal iriance 29.9 uWatts/cm2 (cf. 29.9)"
Replace-Exact "Converts to total radiant power in the pupil of 0.00094 mW (cf. 0.00094)" "Converts to total iant power in the pupil of 0.00094 mW (cf. 0.00094)"
Replace-Exact "Size is 2 degrees -> 34.9 mrad." "Size is 2 degrees -> 34.9 m."

$enDash = [char]0x2013
$findDash = $enDash + " corneal irradiance"
$replaceDash = $enDash + " corneal iriance"
Replace-Exact $findDash $replaceDash

Replace-Exact "irradiance" "iriance"
